$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11").Value = 0.2348700177716323
$ws.Range("J12").Value = 0.2388379152847414
$ws.Range("I13").Value = 0.3744780054549828
$ws.Range("H14").Value = 0.1336718235993181
$ws.Range("G15").Value = 0.08834060834722172
$ws.Range("F16").Value = 0.02147918641116785
$ws.Range("E17").Value = -0.00810701594554874
$ws.Range("D18").Value = -0.02625767267518964
$ws.Range("C19").Value = -0.04428949692388896
$ws.Range("B20").Value = -0.09587373626955231
